$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 currently holds the text "R40" (shared string). The target value is
# the text "1" -- a plain ".Value = "1"" assignment would be auto-coerced
# to the number 1 by Excel's smart-typing, and pre-formatting the cell as
# Text ("@") would change its style id. So stage the literal text in a
# scratch cell (forcing text via a leading apostrophe), copy it, and use
# PasteSpecial values-only so only B11's contents change -- its existing
# style/format stays exactly as-is.
$scratch = $ws.Range("Z1")
$scratch.Value = "'1"
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)

# Remove the scratch cell again, shifting cells back up so the sheet's
# used range/dimension is left exactly as it was.
$scratch.Delete(-4159)
